$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header row (row 1) labels from Chinese to their English
# equivalents. Row 2 (the `{.xxx}` fill template placeholders) is unchanged.
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "String 1"
$ws.Range("D1").Value = "String 2"
$ws.Range("E1").Value = "Image"

# B1 previously had no explicit cell style; align it with the rest of the
# header row by copying the cell format from a neighboring, already-styled
# cell (xlPasteFormats = -4122).
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection moves from E1:E2 (active cell E2) to just E1.
[void]$ws.Range("E1").Select()
